{"js": "// Update the date line and the 25 \"dividend\u00f7divisor=quotient, remainder\"\n// answer cells in the table. Every source text below is unique in the\n// document, so a plain search()-based replace is unambiguous.\nconst replacements = [\n  [\"2025-10-15 Wednesday\", \"2025-10-16 Thursday\"],\n  [\"387\u00f73=129, 0\", \"571\u00f72=285, 1\"],\n  [\"829\u00f79=92, 1\", \"385\u00f76=64, 1\"],\n  [\"729\u00f77=104, 1\", \"544\u00f75=108, 4\"],\n  [\"563\u00f74=140, 3\", \"516\u00f79=57, 3\"],\n  [\"906\u00f76=151, 0\", \"702\u00f78=87, 6\"],\n  [\"359\u00f76=59, 5\", \"574\u00f74=143, 2\"],\n  [\"293\u00f74=73, 1\", \"350\u00f78=43, 6\"],\n  [\"237\u00f72=118, 1\", \"993\u00f72=496, 1\"],\n  [\"415\u00f76=69, 1\", \"900\u00f74=225, 0\"],\n  [\"433\u00f77=61, 6\", \"578\u00f78=72, 2\"],\n  [\"503\u00f76=83, 5\", \"429\u00f72=214, 1\"],\n  [\"335\u00f73=111, 2\", \"874\u00f72=437, 0\"],\n  [\"290\u00f72=145, 0\", \"142\u00f78=17, 6\"],\n  [\"533\u00f79=59, 2\", \"310\u00f76=51, 4\"],\n  [\"741\u00f73=247, 0\", \"646\u00f76=107, 4\"],\n  [\"600\u00f78=75, 0\", \"509\u00f74=127, 1\"],\n  [\"267\u00f78=33, 3\", \"196\u00f74=49, 0\"],\n  [\"755\u00f78=94, 3\", \"555\u00f72=277, 1\"],\n  [\"613\u00f75=122, 3\", \"471\u00f78=58, 7\"],\n  [\"530\u00f73=176, 2\", \"542\u00f78=67, 6\"],\n  [\"521\u00f75=104, 1\", \"648\u00f72=324, 0\"],\n  [\"118\u00f75=23, 3\", \"890\u00f73=296, 2\"],\n  [\"856\u00f72=428, 0\", \"128\u00f76=21, 2\"],\n  [\"449\u00f72=224, 1\", \"732\u00f77=104, 4\"],\n  [\"113\u00f79=12, 5\", \"588\u00f77=84, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 \"dividend\u00f7divisor=quotient, remainder\"\n# answer cells in the table. Every source text below is unique in the\n# document, so a plain Find/Replace is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-10-15 Wednesday\", \"2025-10-16 Thursday\"),\n    @(\"387\u00f73=129, 0\", \"571\u00f72=285, 1\"),\n    @(\"829\u00f79=92, 1\", \"385\u00f76=64, 1\"),\n    @(\"729\u00f77=104, 1\", \"544\u00f75=108, 4\"),\n    @(\"563\u00f74=140, 3\", \"516\u00f79=57, 3\"),\n    @(\"906\u00f76=151, 0\", \"702\u00f78=87, 6\"),\n    @(\"359\u00f76=59, 5\", \"574\u00f74=143, 2\"),\n    @(\"293\u00f74=73, 1\", \"350\u00f78=43, 6\"),\n    @(\"237\u00f72=118, 1\", \"993\u00f72=496, 1\"),\n    @(\"415\u00f76=69, 1\", \"900\u00f74=225, 0\"),\n    @(\"433\u00f77=61, 6\", \"578\u00f78=72, 2\"),\n    @(\"503\u00f76=83, 5\", \"429\u00f72=214, 1\"),\n    @(\"335\u00f73=111, 2\", \"874\u00f72=437, 0\"),\n    @(\"290\u00f72=145, 0\", \"142\u00f78=17, 6\"),\n    @(\"533\u00f79=59, 2\", \"310\u00f76=51, 4\"),\n    @(\"741\u00f73=247, 0\", \"646\u00f76=107, 4\"),\n    @(\"600\u00f78=75, 0\", \"509\u00f74=127, 1\"),\n    @(\"267\u00f78=33, 3\", \"196\u00f74=49, 0\"),\n    @(\"755\u00f78=94, 3\", \"555\u00f72=277, 1\"),\n    @(\"613\u00f75=122, 3\", \"471\u00f78=58, 7\"),\n    @(\"530\u00f73=176, 2\", \"542\u00f78=67, 6\"),\n    @(\"521\u00f75=104, 1\", \"648\u00f72=324, 0\"),\n    @(\"118\u00f75=23, 3\", \"890\u00f73=296, 2\"),\n    @(\"856\u00f72=428, 0\", \"128\u00f76=21, 2\"),\n    @(\"449\u00f72=224, 1\", \"732\u00f77=104, 4\"),\n    @(\"113\u00f79=12, 5\", \"588\u00f77=84, 0\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $oldText,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $newText,\n        2\n    )\n}\n"}
